# Generate Report for Handback
#
# Fills in the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns for each localized file on the
# per-locale sheets, flips the status text from "Ready for handoff" to
# "Handed back: in sync with en-US" everywhere it appears, and widens the
# columns that now hold the longer text / new hyperlinks.

$wb = $excel.ActiveWorkbook

$githubBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/990aedbfbb28a4a6e117e4037a60e4625d824cfc/e2e/"

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: update the per-locale status cells (zh-cn / de-de cols)
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value2 = $newStatus
$overview.Range("F2").Value2 = $newStatus
$overview.Range("E3").Value2 = $newStatus
$overview.Range("F3").Value2 = $newStatus

$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527

# ---------------------------------------------------------------------
# Per-locale sheets: zh-cn + de-de, each with the same two data rows
# (row 2 = 3994de1f..., row 3 = 5ffd0142...) and same handback info,
# except the handoff-xlf name / handback datetime differ by locale.
# ---------------------------------------------------------------------
$locales = @(
    @{
        Sheet = "zh-cn"
        Row2TargetXlf = "3994de1f-0cda-4f3f-b2b8-f0692ebc417e.d7120d0f375ef3487bf647047190f52a22330932.zh-cn.xlf"
        Row3TargetXlf = "5ffd0142-ccb6-4c91-9c0e-c04c855f5fef.09b1fd82aade9b6dd7ae4881b06285504128066e.zh-cn.xlf"
        HandbackDateTime = "2016-08-18 06:45:49"
    },
    @{
        Sheet = "de-de"
        Row2TargetXlf = "3994de1f-0cda-4f3f-b2b8-f0692ebc417e.d7120d0f375ef3487bf647047190f52a22330932.de-de.xlf"
        Row3TargetXlf = "5ffd0142-ccb6-4c91-9c0e-c04c855f5fef.09b1fd82aade9b6dd7ae4881b06285504128066e.de-de.xlf"
        HandbackDateTime = "2016-08-18 06:45:56"
    }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Sheet)

    # Status column (C) -> handed back
    $ws.Range("C2").Value2 = $newStatus
    $ws.Range("C3").Value2 = $newStatus

    # Row 2 -> 3994de1f-...md
    $ws.Hyperlinks.Add($ws.Range("I2"), ($githubBase + "3994de1f-0cda-4f3f-b2b8-f0692ebc417e.md"), "", "", "3994de1f-0cda-4f3f-b2b8-f0692ebc417e.md")
    $ws.Range("J2").Value2 = $locale.Row2TargetXlf
    $ws.Range("K2").Value2 = $locale.HandbackDateTime

    # Row 3 -> 5ffd0142-...md
    $ws.Hyperlinks.Add($ws.Range("I3"), ($githubBase + "5ffd0142-ccb6-4c91-9c0e-c04c855f5fef.md"), "", "", "5ffd0142-ccb6-4c91-9c0e-c04c855f5fef.md")
    $ws.Range("J3").Value2 = $locale.Row3TargetXlf
    $ws.Range("K3").Value2 = $locale.HandbackDateTime

    # Widen columns: Status (C) now holds longer text; Latest Target File (I)
    # and Latest Handback File (J) now hold full file names / hyperlinks.
    $ws.Columns.Item(3).ColumnWidth = 29.9777047293527
    $ws.Columns.Item(9).ColumnWidth = 40
    $ws.Columns.Item(10).ColumnWidth = 40
}
